$wb = $excel.ActiveWorkbook

# --- Sheet: Step1_Data ---
$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("D2").Value2 = 0.009827575114579414
$ws.Range("E2").Value2 = 0.1889429192674925
$ws.Range("G2").Value2 = 0.2049272198497094
$ws.Range("H2").Value2 = 0.1264212815976622
$ws.Range("I2").Value2 = 0.01491137815025717
$ws.Range("J2").Value2 = 0.03261463869147573
$ws.Range("K2").Value2 = 0.01356026304915411
$ws.Range("O2").Value2 = 0.05390148670468743
$ws.Range("P2").Value2 = 0.1535209951083063
$ws.Range("Q2").Value2 = 0.04973915105478471
$ws.Range("R2").Value2 = 0.007705689212730363
$ws.Range("S2").Value2 = 0.003426564452126708
$ws.Range("V2").Value2 = 0.03633486132684614
$ws.Range("W2").Value2 = 0.04700804420976614
$ws.Range("Z2").Value2 = 0.01363392504535475
$ws.Range("AB2").Value2 = 0.02518008097900139
$ws.Range("AD2").Value2 = 0.01410595597829507
$ws.Range("AI2").Value2 = 0.004237970207770404
$ws.Range("E3").Value2 = 0.2813513430161754
$ws.Range("G3").Value2 = 0.1131555183395425
$ws.Range("H3").Value2 = 0.05205739151757476
$ws.Range("K3").Value2 = 0.05197934338320705
$ws.Range("L3").Value2 = 0.01946400501331369
$ws.Range("O3").Value2 = 0.07550251761402933
$ws.Range("P3").Value2 = 0.1014494699892957
$ws.Range("S3").Value2 = 0.009162291893857076
$ws.Range("T3").Value2 = 0.01625408256139031
$ws.Range("V3").Value2 = 0.05640431581689827
$ws.Range("W3").Value2 = 0.01961720254675339
$ws.Range("Z3").Value2 = 0.03153475142411076
$ws.Range("AA3").Value2 = 0.02072466647527125
$ws.Range("AB3").Value2 = 0.05024717635098329
$ws.Range("AD3").Value2 = 0.07304198947904901
$ws.Range("AF3").Value2 = 0.02805393457854805
$ws.Range("D4").Value2 = 0.02848553347204516
$ws.Range("E4").Value2 = 0.1392792210493116
$ws.Range("G4").Value2 = 0.155673735402329
$ws.Range("H4").Value2 = 0.2708867255471549
$ws.Range("J4").Value2 = 0.02868304800364113
$ws.Range("O4").Value2 = 0.06499928558849774
$ws.Range("P4").Value2 = 0.130272359328589
$ws.Range("Q4").Value2 = 0.06955732305727774
$ws.Range("R4").Value2 = 0.01215088489782958
$ws.Range("S4").Value2 = 0.001351986037636072
$ws.Range("U4").Value2 = 0.01131263921515008
$ws.Range("V4").Value2 = 0.0340386388737161
$ws.Range("W4").Value2 = 0.02158803438510002
$ws.Range("Z4").Value2 = 0.01716447513663188
$ws.Range("AD4").Value2 = 0.007436841066311446
$ws.Range("AE4").Value2 = 0.00132501524231182
$ws.Range("AI4").Value2 = 0.005794253696466875
$ws.Range("D5").Value2 = 0.2704086414795412
$ws.Range("E5").Value2 = 0.03111892183841836
$ws.Range("F5").Value2 = 0.2082307508877682
$ws.Range("G5").Value2 = 0.1517049537056725
$ws.Range("I5").Value2 = 0.007756116428900603
$ws.Range("J5").Value2 = 0.006220979023596916
$ws.Range("K5").Value2 = 0.01135064283970488
$ws.Range("M5").Value2 = 0.045846610700101
$ws.Range("N5").Value2 = 0.07565962763983446
$ws.Range("O5").Value2 = 0.007409875756409986
$ws.Range("Q5").Value2 = 0.1074547215308998
$ws.Range("R5").Value2 = 0.06121159919154916
$ws.Range("T5").Value2 = 0.001389977961746261
$ws.Range("W5").Value2 = 0.003286276469417352
$ws.Range("AB5").Value2 = 0.004168017652587896
$ws.Range("AC5").Value2 = 0.006782286893851347
$ws.Range("D6").Value2 = 0.004615369903935967
$ws.Range("E6").Value2 = 0.1958121936262147
$ws.Range("F6").Value2 = 0.0142318632788916
$ws.Range("G6").Value2 = 0.2136483038440998
$ws.Range("H6").Value2 = 0.1586677175086092
$ws.Range("J6").Value2 = 0.01504855555188794
$ws.Range("L6").Value2 = 0.04784183368738099
$ws.Range("N6").Value2 = 0.05963179933118884
$ws.Range("O6").Value2 = 0.1115233741026813
$ws.Range("Q6").Value2 = 0.000820969686727025
$ws.Range("R6").Value2 = 0.05376202541694088
$ws.Range("S6").Value2 = 0.06113406867383525
$ws.Range("U6").Value2 = 0.04018219752908593
$ws.Range("V6").Value2 = 0.007095417975370935
$ws.Range("AC6").Value2 = 0.0159843098831496

# --- Sheet: Step2_Sj ---
$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("D2").Value2 = 0.009827575114579414
$ws.Range("E2").Value2 = 0.1987704943820719
$ws.Range("F2").Value2 = 0.1987704943820719
$ws.Range("G2").Value2 = 0.4036977142317813
$ws.Range("H2").Value2 = 0.5301189958294434
$ws.Range("I2").Value2 = 0.5450303739797007
$ws.Range("J2").Value2 = 0.5776450126711764
$ws.Range("K2").Value2 = 0.5912052757203305
$ws.Range("L2").Value2 = 0.5912052757203305
$ws.Range("M2").Value2 = 0.5912052757203305
$ws.Range("N2").Value2 = 0.5912052757203305
$ws.Range("O2").Value2 = 0.645106762425018
$ws.Range("P2").Value2 = 0.7986277575333243
$ws.Range("Q2").Value2 = 0.8483669085881089
$ws.Range("R2").Value2 = 0.8560725978008393
$ws.Range("S2").Value2 = 0.8594991622529661
$ws.Range("T2").Value2 = 0.8594991622529661
$ws.Range("U2").Value2 = 0.8594991622529661
$ws.Range("V2").Value2 = 0.8958340235798122
$ws.Range("W2").Value2 = 0.9428420677895784
$ws.Range("X2").Value2 = 0.9428420677895784
$ws.Range("Y2").Value2 = 0.9428420677895784
$ws.Range("Z2").Value2 = 0.9564759928349331
$ws.Range("AA2").Value2 = 0.9564759928349331
$ws.Range("AB2").Value2 = 0.9816560738139345
$ws.Range("AC2").Value2 = 0.9816560738139345
$ws.Range("AD2").Value2 = 0.9957620297922296
$ws.Range("AE2").Value2 = 0.9957620297922296
$ws.Range("AF2").Value2 = 0.9957620297922296
$ws.Range("AG2").Value2 = 0.9957620297922296
$ws.Range("AH2").Value2 = 0.9957620297922296
$ws.Range("E3").Value2 = 0.2813513430161754
$ws.Range("F3").Value2 = 0.2813513430161754
$ws.Range("G3").Value2 = 0.39450686135571794
$ws.Range("H3").Value2 = 0.4465642528732927
$ws.Range("I3").Value2 = 0.4465642528732927
$ws.Range("J3").Value2 = 0.4465642528732927
$ws.Range("K3").Value2 = 0.49854359625649974
$ws.Range("L3").Value2 = 0.5180076012698134
$ws.Range("M3").Value2 = 0.5180076012698134
$ws.Range("N3").Value2 = 0.5180076012698134
$ws.Range("O3").Value2 = 0.5935101188838428
$ws.Range("P3").Value2 = 0.6949595888731385
$ws.Range("Q3").Value2 = 0.6949595888731385
$ws.Range("R3").Value2 = 0.6949595888731385
$ws.Range("S3").Value2 = 0.7041218807669956
$ws.Range("T3").Value2 = 0.7203759633283858
$ws.Range("U3").Value2 = 0.7203759633283858
$ws.Range("V3").Value2 = 0.7767802791452841
$ws.Range("W3").Value2 = 0.7963974816920375
$ws.Range("X3").Value2 = 0.7963974816920375
$ws.Range("Y3").Value2 = 0.7963974816920375
$ws.Range("Z3").Value2 = 0.8279322331161483
$ws.Range("AA3").Value2 = 0.8486568995914195
$ws.Range("AB3").Value2 = 0.8989040759424028
$ws.Range("AC3").Value2 = 0.8989040759424028
$ws.Range("AD3").Value2 = 0.9719460654214519
$ws.Range("AE3").Value2 = 0.9719460654214519
$ws.Range("D4").Value2 = 0.02848553347204516
$ws.Range("E4").Value2 = 0.16776475452135675
$ws.Range("F4").Value2 = 0.16776475452135675
$ws.Range("G4").Value2 = 0.32343848992368573
$ws.Range("H4").Value2 = 0.5943252154708407
$ws.Range("I4").Value2 = 0.5943252154708407
$ws.Range("J4").Value2 = 0.6230082634744818
$ws.Range("K4").Value2 = 0.6230082634744818
$ws.Range("L4").Value2 = 0.6230082634744818
$ws.Range("M4").Value2 = 0.6230082634744818
$ws.Range("N4").Value2 = 0.6230082634744818
$ws.Range("O4").Value2 = 0.6880075490629796
$ws.Range("P4").Value2 = 0.8182799083915686
$ws.Range("Q4").Value2 = 0.8878372314488463
$ws.Range("R4").Value2 = 0.8999881163466759
$ws.Range("S4").Value2 = 0.901340102384312
$ws.Range("T4").Value2 = 0.901340102384312
$ws.Range("U4").Value2 = 0.9126527415994621
$ws.Range("V4").Value2 = 0.9466913804731781
$ws.Range("W4").Value2 = 0.9682794148582782
$ws.Range("X4").Value2 = 0.9682794148582782
$ws.Range("Y4").Value2 = 0.9682794148582782
$ws.Range("Z4").Value2 = 0.9854438899949101
$ws.Range("AA4").Value2 = 0.9854438899949101
$ws.Range("AB4").Value2 = 0.9854438899949101
$ws.Range("AC4").Value2 = 0.9854438899949101
$ws.Range("AD4").Value2 = 0.9928807310612215
$ws.Range("AE4").Value2 = 0.9942057463035333
$ws.Range("AF4").Value2 = 0.9942057463035333
$ws.Range("AG4").Value2 = 0.9942057463035333
$ws.Range("AH4").Value2 = 0.9942057463035333
$ws.Range("D5").Value2 = 0.2704086414795412
$ws.Range("E5").Value2 = 0.30152756331795955
$ws.Range("F5").Value2 = 0.5097583142057278
$ws.Range("G5").Value2 = 0.6614632679114003
$ws.Range("H5").Value2 = 0.6614632679114003
$ws.Range("I5").Value2 = 0.6692193843403009
$ws.Range("J5").Value2 = 0.6754403633638978
$ws.Range("K5").Value2 = 0.6867910062036027
$ws.Range("L5").Value2 = 0.6867910062036027
$ws.Range("M5").Value2 = 0.7326376169037037
$ws.Range("N5").Value2 = 0.8082972445435381
$ws.Range("O5").Value2 = 0.8157071202999481
$ws.Range("P5").Value2 = 0.8157071202999481
$ws.Range("Q5").Value2 = 0.9231618418308478
$ws.Range("R5").Value2 = 0.984373441022397
$ws.Range("S5").Value2 = 0.984373441022397
$ws.Range("T5").Value2 = 0.9857634189841432
$ws.Range("U5").Value2 = 0.9857634189841432
$ws.Range("V5").Value2 = 0.9857634189841432
$ws.Range("W5").Value2 = 0.9890496954535606
$ws.Range("X5").Value2 = 0.9890496954535606
$ws.Range("Y5").Value2 = 0.9890496954535606
$ws.Range("Z5").Value2 = 0.9890496954535606
$ws.Range("AA5").Value2 = 0.9890496954535606
$ws.Range("AB5").Value2 = 0.9932177131061485
$ws.Range("D6").Value2 = 0.004615369903935967
$ws.Range("E6").Value2 = 0.20042756353015065
$ws.Range("F6").Value2 = 0.21465942680904224
$ws.Range("G6").Value2 = 0.42830773065314204
$ws.Range("H6").Value2 = 0.5869754481617513
$ws.Range("I6").Value2 = 0.5869754481617513
$ws.Range("J6").Value2 = 0.6020240037136392
$ws.Range("K6").Value2 = 0.6020240037136392
$ws.Range("L6").Value2 = 0.6498658374010202
$ws.Range("M6").Value2 = 0.6498658374010202
$ws.Range("N6").Value2 = 0.709497636732209
$ws.Range("O6").Value2 = 0.8210210108348903
$ws.Range("P6").Value2 = 0.8210210108348903
$ws.Range("Q6").Value2 = 0.8218419805216174
$ws.Range("R6").Value2 = 0.8756040059385582
$ws.Range("S6").Value2 = 0.9367380746123934
$ws.Range("T6").Value2 = 0.9367380746123934
$ws.Range("U6").Value2 = 0.9769202721414794
$ws.Range("V6").Value2 = 0.9840156901168503
$ws.Range("W6").Value2 = 0.9840156901168503
$ws.Range("X6").Value2 = 0.9840156901168503
$ws.Range("Y6").Value2 = 0.9840156901168503
$ws.Range("Z6").Value2 = 0.9840156901168503
$ws.Range("AA6").Value2 = 0.9840156901168503
$ws.Range("AB6").Value2 = 0.9840156901168503

# --- Sheet: Step3_DataPts_0.5 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("D2").Value2 = 7
$ws.Range("F2").Value2 = 0.5301189958294434
$ws.Range("G2").Value2 = 5
$ws.Range("D3").Value2 = 11
$ws.Range("F3").Value2 = 0.5180076012698134
$ws.Range("G3").Value2 = 9
$ws.Range("F4").Value2 = 0.5943252154708407
$ws.Range("F5").Value2 = 0.5097583142057278
$ws.Range("D6").Value2 = 7
$ws.Range("F6").Value2 = 0.5869754481617513
$ws.Range("G6").Value2 = 5

# --- Sheet: Step3_DataPts_0.7 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("D2").Value2 = 15
$ws.Range("F2").Value2 = 0.7986277575333243
$ws.Range("G2").Value2 = 13
$ws.Range("D3").Value2 = 18
$ws.Range("F3").Value2 = 0.7041218807669956
$ws.Range("G3").Value2 = 16
$ws.Range("D4").Value2 = 15
$ws.Range("F4").Value2 = 0.8182799083915686
$ws.Range("G4").Value2 = 13
$ws.Range("D5").Value2 = 12
$ws.Range("F5").Value2 = 0.7326376169037037
$ws.Range("G5").Value2 = 11
$ws.Range("D6").Value2 = 13
$ws.Range("F6").Value2 = 0.709497636732209
$ws.Range("G6").Value2 = 11

# --- Sheet: Step3_DataPts_0.8 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("D2").Value2 = 16
$ws.Range("F2").Value2 = 0.8483669085881089
$ws.Range("G2").Value2 = 14
$ws.Range("D3").Value2 = 25
$ws.Range("F3").Value2 = 0.8279322331161483
$ws.Range("G3").Value2 = 23
$ws.Range("F4").Value2 = 0.8182799083915686
$ws.Range("F5").Value2 = 0.8082972445435381
$ws.Range("F6").Value2 = 0.8210210108348903

# --- Sheet: Step3_DataPts_0.9 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("D2").Value2 = 22
$ws.Range("F2").Value2 = 0.9428420677895784
$ws.Range("G2").Value2 = 20
$ws.Range("D3").Value2 = 29
$ws.Range("F3").Value2 = 0.9719460654214519
$ws.Range("G3").Value2 = 27
$ws.Range("D4").Value2 = 18
$ws.Range("F4").Value2 = 0.901340102384312
$ws.Range("G4").Value2 = 16
$ws.Range("F5").Value2 = 0.9231618418308478
$ws.Range("D6").Value2 = 18
$ws.Range("F6").Value2 = 0.9367380746123934
$ws.Range("G6").Value2 = 16
